$wb = $excel.ActiveWorkbook

# --- ip_address_list (sheet1) ---
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Range("A1:E12").ClearContents()
$ws1.Range("A1:E12").NumberFormat = "@"

$ws1.Range("A1").Value = "529_Witte"
$ws1.Range("B1").Value = "192.168.0.240"
$ws1.Range("C1").Value = "255.255.255.0"
$ws1.Range("D1").Value = "Kamera VS-S160MX :192.168.0.18"
$ws1.Range("E1").Value = "1"

$ws1.Range("A2").Value = "440_Austin"
$ws1.Range("B2").Value = "10.96.205.240"
$ws1.Range("C2").Value = "255.255.255.0"
$ws1.Range("D2").Value = "FortiClient Austin: `npass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK`nFH-2050-20`n10.96.205.80"
$ws1.Range("E2").Value = "0"

$ws1.Range("A3").Value = "474 B_Austin"
$ws1.Range("B3").Value = "10.96.205.175"
$ws1.Range("C3").Value = "255.255.255.0"
$ws1.Range("D3").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.245`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws1.Range("E3").Value = "0"

$ws1.Range("A4").Value = "497_Edcha"
$ws1.Range("B4").Value = "172.26.7.240"
$ws1.Range("C4").Value = "255.255.255.0"
$ws1.Range("D4").Value = "FortiClient Edcha Ex2p78kxp30"
$ws1.Range("E4").Value = "0"

$ws1.Range("A5").Value = "503_Witte"
$ws1.Range("B5").Value = "192.168.0.240"
$ws1.Range("C5").Value = "255.255.255.0"
$ws1.Range("D5").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.267`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws1.Range("E5").Value = "0"

$ws1.Range("A6").Value = "511_Teleflex"
$ws1.Range("B6").Value = "192.168.1.242"
$ws1.Range("C6").Value = "255.255.255.0"
$ws1.Range("D6").Value = "Teleflex d"
$ws1.Range("E6").Value = "0"

$ws1.Range("A7").Value = "514_Teleflex"
$ws1.Range("B7").Value = "192.168.14.240"
$ws1.Range("C7").Value = "255.255.255.0"
$ws1.Range("D7").Value = "PC:192.168.14.240`nCAM: 192.168.14.??NAS:192.168.14.245`n*******************************`nuser: Vision`npass: *Jhv2708"
$ws1.Range("E7").Value = "1"

$ws1.Range("A8").Value = "515_ZF Stara kkkBoleslav"
$ws1.Range("B8").Value = "10.9.250.240"
$ws1.Range("C8").Value = "255.255.255.0"
$ws1.Range("D8").Value = "NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass`n123TPV456"
$ws1.Range("E8").Value = "0"

$ws1.Range("A9").Value = "518_Valeo"
$ws1.Range("B9").Value = "192.168.208.242"
$ws1.Range("C9").Value = "255.255.255.0"
$ws1.Range("E9").Value = "0"

$ws1.Range("A10").Value = "518_Valeo II"
$ws1.Range("B10").Value = "192.168.1.243"
$ws1.Range("C10").Value = "255.255.255.0"
$ws1.Range("E10").Value = "1"

$ws1.Range("A11").Value = "527_Teijin"
$ws1.Range("B11").Value = "10.101.28.176"
$ws1.Range("C11").Value = "255.255.255.0"
$ws1.Range("D11").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.117"
$ws1.Range("E11").Value = "0"

$ws1.Range("A12").Value = "533valeo"
$ws1.Range("B12").Value = "192.168.227.27"
$ws1.Range("C12").Value = "255.255.255.0"
$ws1.Range("E12").Value = "0"

# --- disk_list (sheet3) ---
$ws3 = $wb.Worksheets.Item("disk_list")
$ws3.Range("A1:F7").ClearContents()
$ws3.Range("A1:F7").NumberFormat = "@"

$ws3.Range("A1").Value = "514_Teleflex"
$ws3.Range("B1").Value = "T"
$ws3.Range("C1").Value = "\\192.168.14.245\Data\Kamery"
$ws3.Range("D1").Value = "Vision"
$ws3.Range("E1").Value = "*Jhv2708"

$ws3.Range("A2").Value = "witte"
$ws3.Range("B2").Value = "W"
$ws3.Range("C2").Value = "\\192.168.0.192\"

$ws3.Range("A3").Value = "Domaci Nas"
$ws3.Range("B3").Value = "S"
$ws3.Range("C3").Value = "\\192.168.1.20\Data"

$ws3.Range("A4").Value = "518_Valeo II"
$ws3.Range("B4").Value = "V"
$ws3.Range("C4").Value = "\\192.168.1.10\10_vision"
$ws3.Range("D4").Value = "jhv_vision"
$ws3.Range("E4").Value = "Jhv*2708"
$ws3.Range("F4").Value = "Druha sít, ixon"

$ws3.Range("A5").Value = "518_Valeo"
$ws3.Range("B5").Value = "V"
$ws3.Range("C5").Value = "\\192.168.208.200\10_vision"
$ws3.Range("D5").Value = "jhv_vision"
$ws3.Range("E5").Value = "Jhv*2708"
$ws3.Range("F5").Value = "první sít, ixon`n\\192.168.208.200\10_vision"

$ws3.Range("A6").Value = "515_ZF"
$ws3.Range("B6").Value = "Z"
$ws3.Range("C6").Value = "\\10.9.250.100\08_Project_ZF_515\kamery"
$ws3.Range("D6").Value = "jhvadmin"
$ws3.Range("E6").Value = "jhvadm1n"

$ws3.Range("A7").Value = "474_B Austin"
$ws3.Range("B7").Value = "P"
$ws3.Range("C7").Value = "\\10.96.205.166\DATA"
$ws3.Range("D7").Value = "jhv_vision"
$ws3.Range("E7").Value = "*Jhv2708"
$ws3.Range("F7").Value = "10.96.205.166`nVisionNas_474B`t`n`t`t`t`t`t`tuser:JHV_Vision, omron `nPass:*Jhv2708"

